$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '46.981.46'
$ws.Range("E2").Value = '  +3.14%  '

# Row 3
$ws.Range("D3").Value = '2.628.75'
$ws.Range("E3").Value = '  +7.71%  '

# Row 4
$ws.Range("E4").Value = '  -0.77%  '

# Row 5
$ws.Range("D5").Value = "'309.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.49%  '

# Row 6
$ws.Range("D6").Value = "'102.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.13%  '

# Row 7
$ws.Range("D7").Value = "'0.607"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.48%  '

# Row 8
$ws.Range("E8").Value = '  -0.34%  '

# Row 9
$ws.Range("D9").Value = "'0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +15.11%  '

# Row 10
$ws.Range("D10").Value = "'40.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +15.89%  '

# Row 11
$ws.Range("D11").Value = "'0.0856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.86%  '

# Row 12
$ws.Range("D12").Value = "'54.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.81%  '

# Row 13
$ws.Range("D13").Value = "'8.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +15.73%  '

# Row 14
$ws.Range("D14").Value = '3.020.66'
$ws.Range("E14").Value = '  +6.60%  '

# Row 15
$ws.Range("E15").Value = '  +2.10%  '

# Row 16
$ws.Range("D16").Value = '2.632.45'
$ws.Range("E16").Value = '  +6.43%  '

# Row 17
$ws.Range("D17").Value = "'0.938"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +11.00%  '

# Row 18
$ws.Range("D18").Value = "'15.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.62%  '

# Row 19
$ws.Range("D19").Value = '47.066.62'
$ws.Range("E19").Value = '  +2.71%  '

# Row 20
$ws.Range("E20").Value = '  +9.58%  '

# Row 21
$ws.Range("D21").Value = "'13.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.13%  '

# Row 22
$ws.Range("D22").Value = "'6.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.76%  '

# Row 23
$ws.Range("D23").Value = "'279.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.31%  '

# Row 24
$ws.Range("D24").Value = "'72.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.78%  '

# Row 25
$ws.Range("D25").Value = "'3.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.33%  '

# Row 26
$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +13.69%  '

# Row 27
$ws.Range("D27").Value = "'29.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +37.25%  '

# Row 28
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.24%  '

# Row 29
$ws.Range("D29").Value = "'4.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.42%  '

# Row 30
$ws.Range("D30").Value = "'10.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.54%  '

# Row 31
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = "'39.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.49%  '

# Row 32
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = "'2.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.93%  '

# Row 33
$ws.Range("E33").Value = '  +17.75%  '

# Row 34
$ws.Range("D34").Value = "'3.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.10%  '

# Row 35
$ws.Range("D35").Value = "'2.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.03%  '

# Row 36
$ws.Range("D36").Value = "'2.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.20%  '

# Row 37
$ws.Range("D37").Value = "'0.0852"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.24%  '

# Row 38
$ws.Range("D38").Value = "'153.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.03%  '

# Row 39
$ws.Range("D39").Value = "'0.124"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.61%  '

# Row 40
$ws.Range("E40").Value = '  +7.58%  '

# Row 41
$ws.Range("D41").Value = "'23.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +44.59%  '

# Row 42
$ws.Range("D42").Value = "'16.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.68%  '

# Row 43
$ws.Range("D43").Value = "'3.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +16.64%  '

# Row 44
$ws.Range("D44").Value = "'0.0336"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.27%  '

# Row 45
$ws.Range("D45").Value = "'4.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.20%  '

# Row 46
$ws.Range("D46").Value = '2.128.21'
$ws.Range("E46").Value = '  +6.11%  '

# Row 47
$ws.Range("D47").Value = "'0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.57%  '

# Row 48
$ws.Range("D48").Value = "'94.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.97%  '

# Row 49
$ws.Range("D49").Value = "'9.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +12.48%  '

# Row 50
$ws.Range("D50").Value = "'1.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.06%  '

# Row 51
$ws.Range("D51").Value = "'110.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.25%  '
